$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.590.53'
$ws.Range('E2').Value = '  +0.49%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.580.12'
$ws.Range('E3').Value = '  +0.51%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.14%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.19'
$ws.Range('E6').Value = '  -1.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3675'
$ws.Range('E7').Value = '  -1.13%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.56'
$ws.Range('E8').Value = '  -2.92%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3363'
$ws.Range('E9').Value = '  -0.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.140'
$ws.Range('E10').Value = '  -0.35%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07485'
$ws.Range('E11').Value = '  -0.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9991'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.06'
$ws.Range('E13').Value = '  -1.22%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.028'
$ws.Range('E14').Value = '  -0.30%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.964'
$ws.Range('E15').Value = '  +0.00%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.579.55'
$ws.Range('E16').Value = '  +0.53%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001119'
$ws.Range('E17').Value = '  -0.40%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '88.78'
$ws.Range('E18').Value = '  -2.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06779'
$ws.Range('E19').Value = '  +0.08%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.425'
$ws.Range('E20').Value = '  +2.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.56'
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.20'
$ws.Range('E23').Value = '  +0.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.588.09'
$ws.Range('E24').Value = '  +0.54%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.386'
$ws.Range('E25').Value = '  +1.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.631'
$ws.Range('E26').Value = '  +0.33%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.66'
$ws.Range('E27').Value = '  +2.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.75'
$ws.Range('E28').Value = '  -1.42%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.028'
$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.78'
$ws.Range('E30').Value = '  -0.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.761.29'
$ws.Range('E31').Value = '  +0.82%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.068'
$ws.Range('E32').Value = '  -0.81%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.229'
$ws.Range('E33').Value = '  -0.32%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.001'
$ws.Range('E34').Value = '  -0.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.836'
$ws.Range('E35').Value = '  +0.55%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08319'
$ws.Range('E36').Value = '  -0.48%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02469'
$ws.Range('E37').Value = '  -0.57%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2278'
$ws.Range('E38').Value = '  -1.20%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.489'
$ws.Range('E39').Value = '  +0.52%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06472'
$ws.Range('E40').Value = '  -0.95%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.308'
$ws.Range('E41').Value = '  -2.50%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.48'
$ws.Range('E42').Value = '  +1.08%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6390'
$ws.Range('E43').Value = '  +2.58%  '

$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.04%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.00'
$ws.Range('E45').Value = '  -0.23%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6226'
$ws.Range('E46').Value = '  +6.05%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.771'
$ws.Range('E47').Value = '  -1.13%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.073'
$ws.Range('E48').Value = '  -0.04%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.59'
$ws.Range('E49').Value = '  -2.97%  '

$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.226'
$ws.Range('E50').Value = '  +0.56%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07267'
$ws.Range('E51').Value = '  -0.88%  '
